$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected/active cell to E4
$ws.Range("E4").Select()

# Fill in the new row of data (row 4) with manuel's info
$ws.Range("A4").Value = "manuel"
$ws.Range("B4").Value = "18:15:14:12S"
$ws.Range("C4").Value = "manuel@email.es"
$ws.Range("D4").Value = "ID5"
$ws.Range("E4").Value = 2

# Add a hyperlink on C4 (email cell) pointing to a mailto link
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:manuel@email.es")

# Keep the same cell style as the other email cell (avoid Excel's auto hyperlink style)
$ws.Range("C4").Style = $ws.Range("C2").Style
